$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.264.20'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '1.898.63'
$ws.Range('E3').Value = '  +0.72%  '
$ws.Range('E4').Value = '  -0.88%  '
$ws.Range('D5').Value = '''315.43'
$ws.Range('E5').Value = '  -0.44%  '
$ws.Range('D6').Value = '''1.004'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('D7').Value = '''0.5136'
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('D8').Value = '''0.3929'
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('D9').Value = '''0.08439'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = '''42.51'
$ws.Range('E10').Value = '  +1.62%  '
$ws.Range('D11').Value = '''1.115'
$ws.Range('E11').Value = '  +0.15%  '
$ws.Range('D12').Value = '''6.248'
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = '1.898.54'
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '''7.318'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').Value = '''1.006'
$ws.Range('E16').Value = '  -0.92%  '
$ws.Range('D17').Value = '''93.37'
$ws.Range('E17').Value = '  +2.11%  '
$ws.Range('E18').Value = '  -0.38%  '
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').Value = '''17.85'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = '''1.005'
$ws.Range('E21').Value = '  -0.73%  '
$ws.Range('D22').Value = '''6.018'
$ws.Range('E22').Value = '  +0.85%  '
$ws.Range('D23').Value = '29.261.92'
$ws.Range('E23').Value = '  +2.27%  '
$ws.Range('D24').Value = '''11.18'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').Value = '2.116.26'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('D27').Value = '''159.51'
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').Value = '''2.446'
$ws.Range('D30').Value = '''128.29'
$ws.Range('E30').Value = '  +1.54%  '
$ws.Range('E31').Value = '  +0.68%  '
$ws.Range('D33').Value = '''6.134'
$ws.Range('E33').Value = '  +5.80%  '
$ws.Range('D34').Value = '''3.660'
$ws.Range('E34').Value = '  +0.84%  '
$ws.Range('D35').Value = '''0.02478'
$ws.Range('E35').Value = '  +1.55%  '
$ws.Range('D36').Value = '''0.06560'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('D37').Value = '''9.048'
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('D38').Value = '''0.2192'
$ws.Range('E38').Value = '  +0.10%  '
$ws.Range('D39').Value = '''1.233'
$ws.Range('E39').Value = '  +3.33%  '
$ws.Range('D40').Value = '''5.123'
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('D41').Value = '''0.6497'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('D42').Value = '''1.232'
$ws.Range('E42').Value = '  -2.72%  '
$ws.Range('E43').Value = '  +0.39%  '
$ws.Range('E44').Value = '  -0.52%  '
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('D46').Value = '''3.676'
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('D47').Value = '''2.050'
$ws.Range('E47').Value = '  +1.93%  '
$ws.Range('E48').Value = '  +1.46%  '
$ws.Range('D49').Value = '''123.31'
$ws.Range('E49').Value = '  +0.49%  '
$ws.Range('E50').Value = '  -2.39%  '
$ws.Range('D51').Value = '''77.65'
$ws.Range('E51').Value = '  +0.57%  '

# Reset style on cells that needed a leading apostrophe to force text storage,
# so the quotePrefix flag does not leave a stray style index behind.
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D12').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D16').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D22').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D42').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D49').Style = "Normal"
$ws.Range('D51').Style = "Normal"
